$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage so that
# numeric-looking strings (e.g. "308.32") are not coerced into real numbers,
# which would change their textual representation (e.g. trailing zeros).
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.NumberFormat = "General"
}

$ws.Range('D2').Value = '44.444.05'
$ws.Range('E2').Value = '  +1.38%  '

$ws.Range('D3').Value = '2.249.35'
$ws.Range('E3').Value = '  +1.18%  '

$ws.Range('E4').Value = '  +0.11%  '

Set-TextValue $ws.Range('D5') '308.32'
$ws.Range('E5').Value = '  +2.19%  '

Set-TextValue $ws.Range('D6') '94.78'
$ws.Range('E6').Value = '  +2.35%  '

Set-TextValue $ws.Range('D7') '0.573'
$ws.Range('E7').Value = '  +1.56%  '

Set-TextValue $ws.Range('D8') '1.01'
$ws.Range('E8').Value = '  +0.18%  '

Set-TextValue $ws.Range('D9') '0.527'
$ws.Range('E9').Value = '  +2.83%  '

Set-TextValue $ws.Range('D10') '35.39'
$ws.Range('E10').Value = '  +5.12%  '

$ws.Range('E11').Value = '  +2.29%  '

Set-TextValue $ws.Range('D12') '7.26'
$ws.Range('E12').Value = '  +3.61%  '

$ws.Range('E13').Value = '  +2.04%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D14') '0.842'
$ws.Range('E14').Value = '  +4.69%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.247.14'
$ws.Range('E15').Value = '  -0.44%  '

Set-TextValue $ws.Range('D16') '13.68'
$ws.Range('E16').Value = '  +2.52%  '

$ws.Range('D17').Value = '44.196.50'
$ws.Range('E17').Value = '  +1.32%  '

$ws.Range('D18').Value = '0.0₃0968'
$ws.Range('E18').Value = '  +2.71%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D19') '6.42'
$ws.Range('E19').Value = '  +5.72%  '

$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D20') '12.26'
$ws.Range('E20').Value = '  +2.36%  '

Set-TextValue $ws.Range('D21') '66.21'
$ws.Range('E21').Value = '  +3.81%  '

Set-TextValue $ws.Range('D22') '3.18'
$ws.Range('E22').Value = '  +10.79%  '

Set-TextValue $ws.Range('D23') '237.71'
$ws.Range('E23').Value = '  +1.69%  '

$ws.Range('E24').Value = '  +6.26%  '

$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('E26').Value = '  +5.50%  '

Set-TextValue $ws.Range('D27') '38.34'
$ws.Range('E27').Value = '  +8.43%  '

Set-TextValue $ws.Range('D28') '9.88'
$ws.Range('E28').Value = '  +2.68%  '

Set-TextValue $ws.Range('D29') '6.00'
$ws.Range('E29').Value = '  +3.64%  '

Set-TextValue $ws.Range('D30') '20.13'
$ws.Range('E30').Value = '  +2.73%  '

Set-TextValue $ws.Range('D31') '153.37'
$ws.Range('E31').Value = '  +2.13%  '

Set-TextValue $ws.Range('D32') '0.0801'
$ws.Range('E32').Value = '  +1.38%  '

$ws.Range('E33').Value = '  +1.27%  '

Set-TextValue $ws.Range('D34') '3.14'
$ws.Range('E34').Value = '  -1.32%  '

$ws.Range('E35').Value = '  +2.93%  '

$ws.Range('E36').Value = '  +4.51%  '

Set-TextValue $ws.Range('D37') '1.81'
$ws.Range('E37').Value = '  +5.34%  '

Set-TextValue $ws.Range('D38') '3.48'
$ws.Range('E38').Value = '  +8.59%  '

Set-TextValue $ws.Range('D39') '14.72'
$ws.Range('E39').Value = '  +3.06%  '

Set-TextValue $ws.Range('D40') '3.85'
$ws.Range('E40').Value = '  +3.82%  '

Set-TextValue $ws.Range('D41') '0.0303'
$ws.Range('E41').Value = '  +3.51%  '

$ws.Range('E42').Value = '  +0.17%  '

$ws.Range('D43').Value = '1.749.00'
$ws.Range('E43').Value = '  +1.57%  '

$ws.Range('E44').Value = '  +6.47%  '

Set-TextValue $ws.Range('D45') '80.85'
$ws.Range('E45').Value = '  -1.97%  '

Set-TextValue $ws.Range('D46') '71.27'
$ws.Range('E46').Value = '  +7.65%  '

Set-TextValue $ws.Range('D47') '100.10'
$ws.Range('E47').Value = '  +1.89%  '

Set-TextValue $ws.Range('D48') '4.93'
$ws.Range('E48').Value = '  +1.31%  '

$ws.Range('E49').Value = '  +9.40%  '

Set-TextValue $ws.Range('D50') '55.80'
$ws.Range('E50').Value = '  +5.41%  '

Set-TextValue $ws.Range('D51') '8.22'
$ws.Range('E51').Value = '  +3.08%  '
